$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A39").Value = "Task 39: Manage project page revamp"

$ws.Range("B38").Copy() | Out-Null
$ws.Range("B39").PasteSpecial(-4122) | Out-Null
$ws.Range("B39").Value = "Complete"

$ws.Range("D38").Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4122) | Out-Null
$ws.Range("D39").Value = 45540

$excel.CutCopyMode = $false
$ws.Range("C42").Select()
